# Regenerate the "K" column (column G) values for each start, replacing the
# previous "Strike#" derived figure with the newly computed K value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 6
    4  = 0
    5  = 6
    6  = 2
    7  = 5
    8  = 5
    9  = 4
    10 = 6
    11 = 7
    12 = 5
    13 = 5
    14 = 7
    15 = 10
    16 = 7
    17 = 6
    18 = 8
    19 = 4
    20 = 7
    21 = 8
    22 = 6
    23 = 4
    24 = 9
    25 = 9
    26 = 9
    27 = 6
    28 = 6
    29 = 14
    30 = 6
    31 = 10
    32 = 8
    33 = 5
    34 = 7
    35 = 7
    36 = 7
    37 = 5
    38 = 2
    39 = 1
    40 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
